# Efna5-Epha4.xlsx -- refresh with new TPM-derived NATMI metrics.
# 1) Rename the "Resolving-Mac" cluster label to "Inflammatory-Mac".
# 2) Push the recomputed ligand/receptor/edge expression statistics that
#    follow from the new TPM values (also reorders which row holds the
#    MuSCs vs. Inflammatory-Mac target-cluster record for three pairs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3227736666666667
$ws.Range("H2").Value = 0.968321
$ws.Range("I2").Value = 0.1416094457286952
$ws.Range("J2").Value = 0.1416094457286952
$ws.Range("M2").Value = 8.540560666666666
$ws.Range("N2").Value = 25.621682
$ws.Range("O2").Value = 0.4159358086620884
$ws.Range("P2").Value = 0.4159358086620884
$ws.Range("Q2").Value = 2.756668081769111
$ws.Range("R2").Value = 24.810012735922
$ws.Range("S2").Value = 0.05890043932335494
$ws.Range("T2").Value = 0.05890043932335494

# Row 3
$ws.Range("G3").Value = 0.3227736666666667
$ws.Range("H3").Value = 0.968321
$ws.Range("I3").Value = 0.1416094457286952
$ws.Range("J3").Value = 0.1416094457286952
$ws.Range("O3").Value = 0.563694901924408
$ws.Range("P3").Value = 0.563694901924408
$ws.Range("Q3").Value = 3.735960481472778
$ws.Range("R3").Value = 33.623644333255
$ws.Range("S3").Value = 0.07982452262160659
$ws.Range("T3").Value = 0.07982452262160659

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.3227736666666667
$ws.Range("H4").Value = 0.968321
$ws.Range("I4").Value = 0.1416094457286952
$ws.Range("J4").Value = 0.1416094457286952
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.026642
$ws.Range("N4").Value = 0.079926
$ws.Range("O4").Value = 0.001297498167494471
$ws.Range("P4").Value = 0.001297498167494471
$ws.Range("Q4").Value = 0.008599336027333333
$ws.Range("R4").Value = 0.07739402424599999
$ws.Range("S4").Value = 0.0001837379963328897
$ws.Range("T4").Value = 0.0001837379963328897

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.3227736666666667
$ws.Range("H5").Value = 0.968321
$ws.Range("I5").Value = 0.1416094457286952
$ws.Range("J5").Value = 0.1416094457286952
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.391608
$ws.Range("N5").Value = 1.174824
$ws.Range("O5").Value = 0.01907179124600912
$ws.Range("P5").Value = 0.01907179124600912
$ws.Range("Q5").Value = 0.126400750056
$ws.Range("R5").Value = 1.137606750504
$ws.Range("S5").Value = 0.002700745787400731
$ws.Range("T5").Value = 0.002700745787400731

# Row 6
$ws.Range("I6").Value = 0.8226066833587575
$ws.Range("J6").Value = 0.8226066833587576
$ws.Range("M6").Value = 8.540560666666666
$ws.Range("N6").Value = 25.621682
$ws.Range("O6").Value = 0.4159358086620884
$ws.Range("P6").Value = 0.4159358086620884
$ws.Range("Q6").Value = 16.01343452900422
$ws.Range("R6").Value = 144.120910761038
$ws.Range("S6").Value = 0.3421515760536633
$ws.Range("T6").Value = 0.3421515760536634

# Row 7
$ws.Range("I7").Value = 0.8226066833587575
$ws.Range("J7").Value = 0.8226066833587576
$ws.Range("O7").Value = 0.563694901924408
$ws.Range("P7").Value = 0.563694901924408
$ws.Range("S7").Value = 0.4636991936982773
$ws.Range("T7").Value = 0.4636991936982774

# Row 8
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.8226066833587575
$ws.Range("J8").Value = 0.8226066833587576
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.026642
$ws.Range("N8").Value = 0.079926
$ws.Range("O8").Value = 0.001297498167494471
$ws.Range("P8").Value = 0.001297498167494471
$ws.Range("Q8").Value = 0.04995338589266667
$ws.Range("R8").Value = 0.449580473034
$ws.Range("S8").Value = 0.001067330664226693
$ws.Range("T8").Value = 0.001067330664226693

# Row 9
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.8226066833587575
$ws.Range("J9").Value = 0.8226066833587576
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.391608
$ws.Range("N9").Value = 1.174824
$ws.Range("O9").Value = 0.01907179124600912
$ws.Range("P9").Value = 0.01907179124600912
$ws.Range("Q9").Value = 0.734259648024
$ws.Range("R9").Value = 6.608336832216001
$ws.Range("S9").Value = 0.01568858294259014
$ws.Range("T9").Value = 0.01568858294259015

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.081563
$ws.Range("H10").Value = 0.244689
$ws.Range("I10").Value = 0.03578387091254728
$ws.Range("J10").Value = 0.03578387091254728
$ws.Range("M10").Value = 8.540560666666666
$ws.Range("N10").Value = 25.621682
$ws.Range("O10").Value = 0.4159358086620884
$ws.Range("P10").Value = 0.4159358086620884
$ws.Range("Q10").Value = 0.6965937496553333
$ws.Range("R10").Value = 6.269343746898
$ws.Range("S10").Value = 0.01488379328507014
$ws.Range("T10").Value = 0.01488379328507014

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.081563
$ws.Range("H11").Value = 0.244689
$ws.Range("I11").Value = 0.03578387091254728
$ws.Range("J11").Value = 0.03578387091254728
$ws.Range("O11").Value = 0.563694901924408
$ws.Range("P11").Value = 0.563694901924408
$ws.Range("Q11").Value = 0.9440551575883332
$ws.Range("R11").Value = 8.496496418294999
$ws.Range("S11").Value = 0.02017118560452401
$ws.Range("T11").Value = 0.02017118560452401

# Row 12
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.081563
$ws.Range("H12").Value = 0.244689
$ws.Range("I12").Value = 0.03578387091254728
$ws.Range("J12").Value = 0.03578387091254728
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.026642
$ws.Range("N12").Value = 0.079926
$ws.Range("O12").Value = 0.001297498167494471
$ws.Range("P12").Value = 0.001297498167494471
$ws.Range("Q12").Value = 0.002173001446
$ws.Range("R12").Value = 0.019557013014
$ws.Range("S12").Value = 0.0000464295069348888
$ws.Range("T12").Value = 0.0000464295069348888

# Row 13
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.081563
$ws.Range("H13").Value = 0.244689
$ws.Range("I13").Value = 0.03578387091254728
$ws.Range("J13").Value = 0.03578387091254728
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.391608
$ws.Range("N13").Value = 1.174824
$ws.Range("O13").Value = 0.01907179124600912
$ws.Range("P13").Value = 0.01907179124600912
$ws.Range("Q13").Value = 0.031940723304
$ws.Range("R13").Value = 0.287466509736
$ws.Range("S13").Value = 0.0006824625160182394
$ws.Range("T13").Value = 0.0006824625160182394
